$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.8722624529955464
$ws.Range("E2").Value = 0.8722624529955464

# Row 3
$ws.Range("D3").Value = 0.2526408452775347
$ws.Range("E3").Value = 0.2526408452775347

# Row 4
$ws.Range("D4").Value = 0.0001484038399906608
$ws.Range("E4").Value = 0.0001484038399906608

# Row 5
$ws.Range("D5").Value = 0.00005014400608294908
$ws.Range("E5").Value = 0.00005014400608294908

# Row 6
$ws.Range("D6").Value = 0.9415323349564944
$ws.Range("E6").Value = 0.9415323349564944

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.0005842360999113275
$ws.Range("E7").Value = 0.9994157639000887

# Row 8
$ws.Range("D8").Value = 0.9999999999748757
$ws.Range("E8").Value = 0.00000000002512434704726729

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.09413175930088773
$ws.Range("E9").Value = 0.9058682406991123

# Row 10
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 0

# Row 11
$ws.Range("D11").Value = 0.9731200233041496
$ws.Range("E11").Value = 0.0268799766958504
$ws.Range("F11").Value = 1.502398133277893
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.9367828023877084
$ws.Range("E12").Value = 0.9367828023877084

# Row 13
$ws.Range("D13").Value = 0.007511059640586045
$ws.Range("E13").Value = 0.007511059640586045

# Row 14
$ws.Range("D14").Value = 0.000001370157685433548
$ws.Range("E14").Value = 0.000001370157685433548

# Row 15
$ws.Range("D15").Value = 0.000001669233506717687
$ws.Range("E15").Value = 0.000001669233506717687

# Row 16
$ws.Range("D16").Value = 0.9204813022381296
$ws.Range("E16").Value = 0.9204813022381296

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.001348624106566025
$ws.Range("E17").Value = 0.998651375893434

# Row 18
$ws.Range("D18").Value = 0.9999999999826328
$ws.Range("E18").Value = 0.00000000001736721877421132

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.04161967433631775
$ws.Range("E19").Value = 0.9583803256636823

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = 0.00000001107088636722107
$ws.Range("E20").Value = 0.9999999889291137

# Row 21
$ws.Range("D21").Value = 0.9922793240571101
$ws.Range("E21").Value = 0.007720675942889943
$ws.Range("F21").Value = 3.341503858566284
$ws.Range("G21").Value = 0.5
